# Trade #198 closed at 2026-02-17 10:08:22 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up numbers for the newly
# closed volatility_scorer trade (#198) and appends the two new trade
# rows (the closed volatility_scorer trade #198 and the freshly opened
# MarketMaking trade #199) to the "All Trades" sheet plus their
# respective per-strategy sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Summary sheet roll-up numbers
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.54    # Current Capital
$summary.Range("B4").Value = -0.46      # Total P&L $
$summary.Range("B6").Value = 198        # Total Trades
$summary.Range("B8").Value = 84         # Losing Trades
$summary.Range("B9").Value = 41.41      # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet - volatility_scorer row (row 12)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C12").Value = 99.22
$status.Range("D12").Value = 17
$status.Range("E12").Value = -0.78
$status.Range("F12").Value = -0.78
$status.Range("G12").Value = 29.41

# ---------------------------------------------------------------------
# 3. Helper to write one trade row given a worksheet + row number
# ---------------------------------------------------------------------
function Set-TradeRow {
    param(
        $Sheet,
        $Row,
        $TradeNum,
        $Date,
        $Time,
        $Strategy,
        $Side,
        $EntryPrice,
        $ExitPrice,
        $Status,
        $PnlPct,
        $PnlDollar,
        $CapitalAfter,
        $EntrySlippage,
        $ExitSlippage,
        $Confidence,
        $EntryReason,
        $ExitReason,
        $Duration
    )

    $Sheet.Cells.Item($Row, 1).Value = $TradeNum
    # Force the Date/Time columns to Text format first so Excel doesn't
    # silently convert "2026-02-17" into a date serial number - the
    # source data models these as plain strings.
    $Sheet.Cells.Item($Row, 2).NumberFormat = "@"
    $Sheet.Cells.Item($Row, 2).Value = $Date
    $Sheet.Cells.Item($Row, 3).NumberFormat = "@"
    $Sheet.Cells.Item($Row, 3).Value = $Time
    $Sheet.Cells.Item($Row, 4).Value = $Strategy
    $Sheet.Cells.Item($Row, 5).Value = $Side
    $Sheet.Cells.Item($Row, 6).Value = $EntryPrice
    if ($null -ne $ExitPrice) {
        $Sheet.Cells.Item($Row, 7).Value = $ExitPrice
    }
    $Sheet.Cells.Item($Row, 8).Value = $Status
    $Sheet.Cells.Item($Row, 9).Value = $PnlPct
    $Sheet.Cells.Item($Row, 10).Value = $PnlDollar
    $Sheet.Cells.Item($Row, 11).Value = $CapitalAfter
    $Sheet.Cells.Item($Row, 12).Value = $EntrySlippage
    $Sheet.Cells.Item($Row, 13).Value = $ExitSlippage
    $Sheet.Cells.Item($Row, 14).Value = $Confidence
    $Sheet.Cells.Item($Row, 15).Value = $EntryReason
    if ($null -ne $ExitReason) {
        $Sheet.Cells.Item($Row, 16).Value = $ExitReason
    }
    $Sheet.Cells.Item($Row, 17).Value = $Duration
}

# ---------------------------------------------------------------------
# 4. All Trades sheet - append trade #198 (row 199) and #199 (row 200)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

Set-TradeRow $allTrades 199 198 `
    "2026-02-17" "10:08:15" "volatility_scorer" "NEUTRAL" `
    0.03 0.02 "CLOSED" `
    -33.3333 -0.01 99.22 `
    0 0 0.85 `
    "Low vol market (score: inf) - ideal for market making" `
    "early_exit" 0.17

Set-TradeRow $allTrades 200 199 `
    "2026-02-17" "10:08:15" "MarketMaking" "UP" `
    0.97 $null "OPEN" `
    0 0 100.3171991854615 `
    0 0 0.6 `
    "Normal spread capture: 19600 bps" `
    $null 0

# ---------------------------------------------------------------------
# 5. volatility_scorer sheet - append trade #198 as local row 18
# ---------------------------------------------------------------------
$volScorer = $wb.Worksheets.Item("volatility_scorer")

Set-TradeRow $volScorer 18 198 `
    "2026-02-17" "10:08:15" "volatility_scorer" "NEUTRAL" `
    0.03 0.02 "CLOSED" `
    -33.3333 -0.01 99.22 `
    0 0 0.85 `
    "Low vol market (score: inf) - ideal for market making" `
    "early_exit" 0.17

# ---------------------------------------------------------------------
# 6. MarketMaking sheet - append trade #199 as local row 183
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

Set-TradeRow $marketMaking 183 199 `
    "2026-02-17" "10:08:15" "MarketMaking" "UP" `
    0.97 $null "OPEN" `
    0 0 100.3171991854615 `
    0 0 0.6 `
    "Normal spread capture: 19600 bps" `
    $null 0
